$wb = $excel.ActiveWorkbook

# "Hoja1" sheet - update the daily conversion summary text in A1
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.08 = 7533.44 pesos`n✅ 7533.44 pesos = 2.08 = 952.79 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# "tasas" sheet - update the input rate cells
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 480
$ws2.Range("O10").Value = 3616.05
$ws2.Range("N12").Value = 3629.97
$ws2.Range("O12").Value = 459.101
